$wb = $excel.ActiveWorkbook

# --- "Logs" sheet: append the new row 12 (new support ticket) ---
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A12").Value = "Retour aanmelden"
$ws.Range("B12").Value = "mailmind.test@zohomail.eu"
$ws.Range("C12").Value = "Ik wil graag een artikel retourneren. Hoe werkt dat?"
$ws.Range("D12").Value = "Retour / Terugbetaling"
$ws.Range("E12").Value = "Beste klant,`nBedankt voor je bericht. Om een artikel te retourneren, volg je eenvoudig de onderstaande stappen:`n1. Log in op je account op onze website.`n2. Ga naar je bestelgeschiedenis en selecteer de bestelling waarvan je een artikel wilt retourneren.`n3. Klik op de optie 'Retourneren' naast het artikel dat je wilt terugsturen.`n4. Volg de instructies om het retourproces te voltooien en het retourlabel te ontvangen.`nZodra we het geretourneerde artikel hebben ontvangen en verwerkt, zullen we het aankoopbedrag terugstorten op de oorspronkelijke betaalmethode.`nMocht je nog vragen hebben of hulp nodig hebben bij het retourneren, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam Bedrijf] E-mailassistent"
$ws.Range("F12").Value = "2025-06-24 19:57:21"
$ws.Range("G12").Value = "Ja"

# Writing the multi-line answer into E12 makes this engine (like real Excel)
# auto-fit the row to the wrapped text, pinning an explicit row height. The
# source workbook doesn't carry one for this row, so re-run AutoFit, which
# clears the "customHeight" flag again (re-measure on open, like the other
# rows) instead of leaving a pinned ht="..." behind.
$ws.Rows.Item(12).AutoFit()

# Conditional formatting ranges need to grow from row 11 to row 12.
$ws.Range("D2:D11").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D12"))
$ws.Range("G2:G11").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G12"))

# --- "Dashboard" sheet: bump the "Retour / Terugbetaling" tally ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 4
